$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) values that look numeric are written as literal text,
# matching the source data which stores prices as inline strings (e.g. "1.004", "83.50").
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "28.927.80"
$ws.Range("E2").Value = "  -2.03%  "
$ws.Range("D3").Value = "1.899.67"
$ws.Range("E3").Value = "  -3.97%  "
$ws.Range("E4").Value = "  -0.26%  "
$ws.Range("D5").Value = "324.32"
$ws.Range("E5").Value = "  -1.03%  "
$ws.Range("E6").Value = "  -0.29%  "
$ws.Range("D7").Value = "0.4585"
$ws.Range("E7").Value = "  -1.68%  "
$ws.Range("D8").Value = "0.3810"
$ws.Range("E8").Value = "  -2.62%  "
$ws.Range("D9").Value = "0.07707"
$ws.Range("E9").Value = "  -3.17%  "
$ws.Range("D10").Value = "0.9748"
$ws.Range("E10").Value = "  -1.97%  "
$ws.Range("D11").Value = "22.01"
$ws.Range("E11").Value = "  -3.87%  "
$ws.Range("D12").Value = "1.885.49"
$ws.Range("E12").Value = "  -4.69%  "
$ws.Range("D13").Value = "6.920"
$ws.Range("E13").Value = "  -3.96%  "
$ws.Range("D14").Value = "5.630"
$ws.Range("E14").Value = "  -3.92%  "
$ws.Range("D15").Value = "0.07033"
$ws.Range("E15").Value = "  -0.77%  "
$ws.Range("D16").Value = "1.003"
$ws.Range("E16").Value = "  -0.38%  "
$ws.Range("D17").Value = "83.50"
$ws.Range("E17").Value = "  -4.89%  "
$ws.Range("D18").Value = "0.000009467"
$ws.Range("E18").Value = "  -4.90%  "
$ws.Range("D19").Value = "16.58"
$ws.Range("E19").Value = "  -4.36%  "
$ws.Range("E20").Value = "  -0.23%  "
$ws.Range("D21").Value = "28.906.88"
$ws.Range("E21").Value = "  -2.09%  "
$ws.Range("D22").Value = "5.281"
$ws.Range("E22").Value = "  -5.07%  "
$ws.Range("E23").Value = "  -3.27%  "
$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").Value = "2.091"
$ws.Range("E24").Value = "  -0.91%  "
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").Value = "157.87"
$ws.Range("E25").Value = "  -0.49%  "
$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").Value = "18.95"
$ws.Range("E26").Value = "  -3.31%  "
$ws.Range("B27").Value = "InternetComputer(DFINITY)"
$ws.Range("C27").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D27").Value = "5.605"
$ws.Range("E27").Value = "  -3.25%  "
$ws.Range("B28").Value = "BitcoinCash"
$ws.Range("C28").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D28").Value = "117.25"
$ws.Range("E28").Value = "  -1.96%  "
$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D29").Value = "1.830"
$ws.Range("E29").Value = "  -4.27%  "
$ws.Range("B30").Value = "Stellar"
$ws.Range("C30").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D30").Value = "0.09232"
$ws.Range("E30").Value = "  -2.12%  "
$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").Value = "0.8580"
$ws.Range("E31").Value = "  -4.31%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "5.074"
$ws.Range("E32").Value = "  -3.20%  "
$ws.Range("B33").Value = "ARBITRUM"
$ws.Range("C33").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D33").Value = "1.236"
$ws.Range("E33").Value = "  -6.63%  "
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "2.988"
$ws.Range("E34").Value = "  -6.43%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "0.05657"
$ws.Range("E35").Value = "  -2.89%  "
$ws.Range("B36").Value = "TrustWalletToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D36").Value = "1.137"
$ws.Range("E36").Value = "  -2.78%  "
$ws.Range("B37").Value = "Frax"
$ws.Range("C37").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D37").Value = "1.002"
$ws.Range("E37").Value = "  -0.06%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "0.02025"
$ws.Range("E38").Value = "  -3.81%  "
$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D39").Value = "0.5466"
$ws.Range("E39").Value = "  -4.55%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "7.365"
$ws.Range("E40").Value = "  -5.55%  "
$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").Value = "0.1747"
$ws.Range("E41").Value = "  -3.05%  "
$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").Value = "9.245"
$ws.Range("E42").Value = "  -4.38%  "
$ws.Range("B43").Value = "MXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D43").Value = "2.753"
$ws.Range("E43").Value = "  -1.49%  "
$ws.Range("B44").Value = "Decentraland"
$ws.Range("C44").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D44").Value = "0.5139"
$ws.Range("E44").Value = "  -4.30%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "11.21"
$ws.Range("E45").Value = "  -5.01%  "
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").Value = "0.06807"
$ws.Range("E46").Value = "  -1.72%  "
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").Value = "2.066"
$ws.Range("E47").Value = "  -5.32%  "
$ws.Range("B48").Value = "PEPE"
$ws.Range("C48").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D48").Value = "0.000002610"
$ws.Range("E48").Value = "  -16.06%  "
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").Value = "109.86"
$ws.Range("E49").Value = "  -3.95%  "
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").Value = "1.765"
$ws.Range("E50").Value = "  -3.66%  "
$ws.Range("B51").Value = "PaxDollar"
$ws.Range("C51").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D51").Value = "1.002"
$ws.Range("E51").Value = "  -0.31%  "

# Restore default (unstyled) formatting on the Price column so no stray
# number-format style is left applied to the cells.
$priceRange.Style = "Normal"
